$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.244.08'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '3.570.98'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''609.30'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('D6').Value = '''145.32'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('D7').Value = '3.569.97'
$ws.Range('E7').Value = '  +2.76%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '''0.493'
$ws.Range('E9').Value = '  +4.05%  '
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').Value = '''7.90'
$ws.Range('E11').Value = '  -3.04%  '
$ws.Range('E12').Value = '  +0.69%  '
$ws.Range('D13').Value = '4.177.26'
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('E14').Value = '  +2.57%  '
$ws.Range('D15').Value = '''29.99'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').Value = '3.554.93'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '66.357.54'
$ws.Range('E17').Value = '  +0.42%  '
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('D19').Value = '''11.49'
$ws.Range('E19').Value = '  +10.72%  '
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').Value = '''14.89'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('D22').Value = '''430.09'
$ws.Range('E22').Value = '  +2.47%  '
$ws.Range('D23').Value = '''0.616'
$ws.Range('E23').Value = '  +4.53%  '
$ws.Range('E24').Value = '  +2.39%  '
$ws.Range('D25').Value = '3.712.21'
$ws.Range('E25').Value = '  +2.62%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +3.90%  '
$ws.Range('E28').Value = '  +2.39%  '
$ws.Range('D29').Value = '''7.95'
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('D30').Value = '''9.10'
$ws.Range('E30').Value = '  -2.33%  '
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('D34').Value = '3.564.24'
$ws.Range('E34').Value = '  +2.68%  '
$ws.Range('E35').Value = '  -5.95%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '''1.74'
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('D38').Value = '''7.88'
$ws.Range('E38').Value = '  +2.42%  '
$ws.Range('E39').Value = '  +0.87%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '''177.76'
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').Value = '''5.23'
$ws.Range('E43').Value = '  +2.72%  '
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').Value = '''1.93'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').Value = '''25.74'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').Value = '''23.55'
$ws.Range('E50').Value = '  +9.19%  '
$ws.Range('D51').Value = '''7.15'
$ws.Range('E51').Value = '  +0.47%  '
